$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.873.95"
$ws.Range("E2").Value = "  +7.06%  "
$ws.Range("D3").Value = "3.628.80"
$ws.Range("E3").Value = "  +6.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.24"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "191.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +7.09%  "
$ws.Range("E7").Value = "  +3.26%  "
$ws.Range("D8").Value = "3.619.47"
$ws.Range("E8").Value = "  +6.22%  "
$ws.Range("E10").Value = "  +2.50%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.47"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +6.99%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000297"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.08%  "
$ws.Range("E14").Value = "  +5.39%  "
$ws.Range("D15").Value = "4.204.82"
$ws.Range("E15").Value = "  +6.14%  "
$ws.Range("D16").Value = "3.631.33"
$ws.Range("E16").Value = "  +5.24%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "19.43"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +6.05%  "
$ws.Range("D18").Value = "70.791.48"
$ws.Range("E18").Value = "  +6.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.81%  "
$ws.Range("E20").Value = "  +0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.06"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.60"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +18.39%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +8.14%  "
$ws.Range("E25").Value = "  +8.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.55%  "
$ws.Range("E27").Value = "  +6.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.26"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.47%  "
$ws.Range("E29").Value = "  +7.00%  "
$ws.Range("E30").Value = "  +3.58%  "
$ws.Range("E31").Value = "  +11.56%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "624.24"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.94%  "
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "12.26"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.94%  "
$ws.Range("E34").Value = "  +8.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "65.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.30%  "
$ws.Range("D36").Value = "0.0₃0832"
$ws.Range("E36").Value = "  +8.28%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.411"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.51%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "38.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +5.06%  "
$ws.Range("E39").Value = "  +1.79%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  +2.52%  "
$ws.Range("D42").Value = "3.328.91"
$ws.Range("E42").Value = "  +6.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.09"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.80%  "
$ws.Range("E44").Value = "  +6.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.36"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.10%  "
$ws.Range("E47").Value = "  +2.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.19"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.36%  "
$ws.Range("E49").Value = "  +3.01%  "
$ws.Range("E50").Value = "  +4.73%  "
$ws.Range("B51").Value = "FirstDigitalUSD"
$ws.Range("C51").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.999"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.06%  "
